$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new changelog entry as a new shared string, in cell F29
# (new week's activity note), matching the style (wrap text) used by
# the other "Activity" cells in column F.
$newEntry = @"
* Added a runtime Dialogue Graph object
* Parsed editor data into runtime data (better suited for using in a project)
* Fixed an issue where the port capacity was not taken into account (i.e if a port could only accept one connection at a time, the tool would ignore that and default to more than one connection)
* Began working on a custom inspector for the runtime monobehaviour
* Improved/modified node styling
"@

$ws.Cells.Item(29, 6).Value = $newEntry
$ws.Cells.Item(29, 6).WrapText = $true
$ws.Rows.Item(29).RowHeight = 15

# Update the timesheet entries for day 29 (row 29: Start/End/Break values)
$ws.Cells.Item(29, 2).Value = 0.40201388888888889
$ws.Cells.Item(29, 3).Value = 0.84842592592592592
$ws.Cells.Item(29, 4).Value = 1 / 72

# Hide the now-empty rows for days 21-28 (no work logged that period)
$ws.Rows.Item(21).Hidden = $true
$ws.Rows.Item(22).Hidden = $true
$ws.Rows.Item(23).Hidden = $true
$ws.Rows.Item(24).Hidden = $true
$ws.Rows.Item(25).Hidden = $true
$ws.Rows.Item(26).Hidden = $true
$ws.Rows.Item(27).Hidden = $true
$ws.Rows.Item(28).Hidden = $true

# Move the sheet selection to match where the author last clicked
$ws.Range("F30").Select()
